# Append a new top entry to the "ランサーズ" sheet (row insert at row 2, shifting
# everything else down by one row) and append a new summary row to "統計".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "ランサーズ" - insert a new case at the top (row 2)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Shift every existing data row (2..27) down to (3..28) and open up row 2 for
# the newly scraped listing.
$ws1.Rows.Item(2).Insert()

$ws1.Range("A2").Value = "2025-08-29 01:45:10"
$ws1.Range("B2").Value = "【急募】フロントエンド Webサービス開発/改修"
$ws1.Range("C2").Value = "システム開発"
$ws1.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws1.Range("E2").Value = "期限情報なし"
$ws1.Range("F2").Value = "https://www.lancers.jp/work/detail/5370186"
$ws1.Range("G2").Value = 68
$ws1.Range("H2").Value = "◆開発"

# The row-insert above shifts cell content but Excel's Hyperlinks collection
# does not renumber its `ref=` anchors automatically, so rebuild the whole
# collection (still in row order -> stable rId1..rId27) once content is final.
$ws1.Hyperlinks.Delete()
for ($r = 2; $r -le 28; $r++) {
    $urlCell = $ws1.Cells.Item($r, 6)
    $url = $urlCell.Value()
    $ws1.Hyperlinks.Add($urlCell, $url)
}
$ws1.Range("F2:F28").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Sheet 2: "統計" - append the matching stats snapshot as a new last row (16)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A16").Value = "2025-08-29T01:45:10.290282"
$ws2.Range("B16").Value = 15
$ws2.Range("C16").Value = "全案件リスト"
$ws2.Range("D16").Value = 53.3
$ws2.Range("E16").Value = 7
$ws2.Range("F16").Value = 4
$ws2.Range("G16").Value = 15
